# A new product ("مجموعه برد") was added to the price list. Because the
# report rows are sorted alphabetically by product name, this new product
# lands between "ماسك جلسات اطفال" (row 118) and "مرطب شفاه لونا جوز هند ابيض"
# (row 119), pushing every row from 119 onward down by one - including the
# totals row and the footer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new data row by inserting a blank row where the totals
# row currently lives (row 124). This shifts the old totals row (124) down
# to 125, and the old footer row (125) down to 126.
$ws.Rows("124:124").Insert()

# Re-create the 3-way merge pattern used by every normal data row
# (name spans B:G, balance spans H:K, price spans L:M).
$ws.Range("B124:G124").Merge()
$ws.Range("H124:K124").Merge()
$ws.Range("L124:M124").Merge()

# Restore the thin light-grey row-separator border (bottom edge) that every
# other data row carries, matching the formatting used throughout the table.
$newRow = $ws.Range("A124:N124")
$newRow.Borders.Item(9).LineStyle = 1
$newRow.Borders.Item(9).Color = 13882323

# Match the auto-fit row heights Excel settles on for the new data row and
# for the (now one-row-taller, re-flowed) totals row beneath it.
$ws.Rows("124:124").RowHeight = 24.75
$ws.Rows("125:125").RowHeight = 26.25

# --- New data row 124: duplicate of the former last row's stats, now
#     holding the item that slid down off the end of the list.
$ws.Range("A124").Value = 121
$ws.Range("B124").Value = "مناديل مبلله كبيره"
$ws.Range("H124").Value = "6:0"
$ws.Range("L124").Value = 30
$ws.Range("N124").Value = "1:0"

# --- Every row from 119 to 123 now shows the product that used to sit one
#     row above it, and the brand-new product takes over row 119.
$ws.Range("B119").Value = "مجموعه برد"
$ws.Range("H119").Value = "-1:0"
$ws.Range("L119").Value = 8

$ws.Range("B120").Value = "مرطب شفاه لونا جوز هند ابيض"
$ws.Range("H120").Value = "1:0"
$ws.Range("L120").Value = 20

$ws.Range("B121").Value = "مرهم راسب ابيض"
$ws.Range("H121").Value = "8:0"
$ws.Range("L121").Value = 8

$ws.Range("B122").Value = "معجون سنسوداين صغير"
$ws.Range("H122").Value = "3:0"
$ws.Range("L122").Value = 40

$ws.Range("B123").Value = "معجون سيجنال 25 مل"
$ws.Range("H123").Value = "4:0"
$ws.Range("L123").Value = 20

# --- Totals row moved from 124 to 125; bump the grand total by the new
#     product's price (7648.1 + 8 = 7656.1).
$ws.Range("K125").Value = 7656.1

Write-Host "Done applying update"
